# Natmi following Dr Hou advice
# Update LR-pair table: re-run with ECs cluster added (M1/M2/Neutro/ECs targets),
# refreshed statistics, and table extended from 9 to 12 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "M1"
$ws.Cells.Item(2,2).Value = "Ccl12"
$ws.Cells.Item(2,3).Value = "Ccr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 40.32940166666666
$ws.Cells.Item(2,8).Value = 120.988205
$ws.Cells.Item(2,9).Value = 0.4705770439863239
$ws.Cells.Item(2,10).Value = 0.4705770439863239
$ws.Cells.Item(2,11).Value = 1.0
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.002632
$ws.Cells.Item(2,14).Value = 0.007896
$ws.Cells.Item(2,15).Value = 0.000006230194780080539
$ws.Cells.Item(2,16).Value = 0.000006230194780080539
$ws.Cells.Item(2,17).Value = 0.1061469851866667
$ws.Cells.Item(2,18).Value = 0.95532286668
$ws.Cells.Item(2,19).Value = 0.000002931786643069325
$ws.Cells.Item(2,20).Value = 0.000002931786643069325

# Row 3
$ws.Cells.Item(3,1).Value = "M1"
$ws.Cells.Item(3,2).Value = "Ccl12"
$ws.Cells.Item(3,3).Value = "Ccr1"
$ws.Cells.Item(3,4).Value = "M1"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 40.32940166666666
$ws.Cells.Item(3,8).Value = 120.988205
$ws.Cells.Item(3,9).Value = 0.4705770439863239
$ws.Cells.Item(3,10).Value = 0.4705770439863239
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 58.61543533333333
$ws.Cells.Item(3,14).Value = 175.846306
$ws.Cells.Item(3,15).Value = 0.1387483203821739
$ws.Cells.Item(3,16).Value = 0.1387483203821739
$ws.Cells.Item(3,17).Value = 2363.925435424525
$ws.Cells.Item(3,18).Value = 21275.32891882073
$ws.Cells.Item(3,19).Value = 0.0652917744635108
$ws.Cells.Item(3,20).Value = 0.06529177446351081

# Row 4
$ws.Cells.Item(4,1).Value = "M1"
$ws.Cells.Item(4,2).Value = "Ccl12"
$ws.Cells.Item(4,3).Value = "Ccr1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 40.32940166666666
$ws.Cells.Item(4,8).Value = 120.988205
$ws.Cells.Item(4,9).Value = 0.4705770439863239
$ws.Cells.Item(4,10).Value = 0.4705770439863239
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 28.58678433333334
$ws.Cells.Item(4,14).Value = 85.76035300000001
$ws.Cells.Item(4,15).Value = 0.0676676422997042
$ws.Cells.Item(4,16).Value = 0.0676676422997042
$ws.Cells.Item(4,17).Value = 1152.887907737374
$ws.Cells.Item(4,18).Value = 10375.99116963637
$ws.Cells.Item(4,19).Value = 0.03184283908691873
$ws.Cells.Item(4,20).Value = 0.03184283908691874

# Row 5
$ws.Cells.Item(5,1).Value = "M1"
$ws.Cells.Item(5,2).Value = "Ccl12"
$ws.Cells.Item(5,3).Value = "Ccr1"
$ws.Cells.Item(5,4).Value = "Neutro"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 40.32940166666666
$ws.Cells.Item(5,8).Value = 120.988205
$ws.Cells.Item(5,9).Value = 0.4705770439863239
$ws.Cells.Item(5,10).Value = 0.4705770439863239
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 335.2538503333333
$ws.Cells.Item(5,14).Value = 1005.761551
$ws.Cells.Item(5,15).Value = 0.7935778071233418
$ws.Cells.Item(5,16).Value = 0.7935778071233417
$ws.Cells.Item(5,17).Value = 13520.58719038955
$ws.Cells.Item(5,18).Value = 121685.284713506
$ws.Cells.Item(5,19).Value = 0.3734394986492512
$ws.Cells.Item(5,20).Value = 0.3734394986492512

# Row 6
$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Ccl12"
$ws.Cells.Item(6,3).Value = "Ccr1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 35.65443533333333
$ws.Cells.Item(6,8).Value = 106.963306
$ws.Cells.Item(6,9).Value = 0.41602796200245
$ws.Cells.Item(6,10).Value = 0.41602796200245
$ws.Cells.Item(6,11).Value = 1.0
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.002632
$ws.Cells.Item(6,14).Value = 0.007896
$ws.Cells.Item(6,15).Value = 0.000006230194780080539
$ws.Cells.Item(6,16).Value = 0.000006230194780080539
$ws.Cells.Item(6,17).Value = 0.09384247379733333
$ws.Cells.Item(6,18).Value = 0.8445822641759999
$ws.Cells.Item(6,19).Value = 0.000002591935237235209
$ws.Cells.Item(6,20).Value = 0.000002591935237235208

# Row 7
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Ccl12"
$ws.Cells.Item(7,3).Value = "Ccr1"
$ws.Cells.Item(7,4).Value = "M1"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 35.65443533333333
$ws.Cells.Item(7,8).Value = 106.963306
$ws.Cells.Item(7,9).Value = 0.41602796200245
$ws.Cells.Item(7,10).Value = 0.41602796200245
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 58.61543533333333
$ws.Cells.Item(7,14).Value = 175.846306
$ws.Cells.Item(7,15).Value = 0.1387483203821739
$ws.Cells.Item(7,16).Value = 0.1387483203821739
$ws.Cells.Item(7,17).Value = 2089.900248627515
$ws.Cells.Item(7,18).Value = 18809.10223764763
$ws.Cells.Item(7,19).Value = 0.05772318095985879
$ws.Cells.Item(7,20).Value = 0.05772318095985879

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Ccl12"
$ws.Cells.Item(8,3).Value = "Ccr1"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 35.65443533333333
$ws.Cells.Item(8,8).Value = 106.963306
$ws.Cells.Item(8,9).Value = 0.41602796200245
$ws.Cells.Item(8,10).Value = 0.41602796200245
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 28.58678433333334
$ws.Cells.Item(8,14).Value = 85.76035300000001
$ws.Cells.Item(8,15).Value = 0.0676676422997042
$ws.Cells.Item(8,16).Value = 0.0676676422997042
$ws.Cells.Item(8,17).Value = 1019.24565340078
$ws.Cells.Item(8,18).Value = 9173.210880607017
$ws.Cells.Item(8,19).Value = 0.02815163131945672
$ws.Cells.Item(8,20).Value = 0.02815163131945672

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Ccl12"
$ws.Cells.Item(9,3).Value = "Ccr1"
$ws.Cells.Item(9,4).Value = "Neutro"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 35.65443533333333
$ws.Cells.Item(9,8).Value = 106.963306
$ws.Cells.Item(9,9).Value = 0.41602796200245
$ws.Cells.Item(9,10).Value = 0.41602796200245
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 335.2538503333333
$ws.Cells.Item(9,14).Value = 1005.761551
$ws.Cells.Item(9,15).Value = 0.7935778071233418
$ws.Cells.Item(9,16).Value = 0.7935778071233417
$ws.Cells.Item(9,17).Value = 11953.28672696085
$ws.Cells.Item(9,18).Value = 107579.5805426476
$ws.Cells.Item(9,19).Value = 0.3301505577878972
$ws.Cells.Item(9,20).Value = 0.3301505577878971

# Row 10
$ws.Cells.Item(10,1).Value = "Neutro"
$ws.Cells.Item(10,2).Value = "Ccl12"
$ws.Cells.Item(10,3).Value = "Ccr1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 9.718179666666666
$ws.Cells.Item(10,8).Value = 29.154539
$ws.Cells.Item(10,9).Value = 0.1133949940112261
$ws.Cells.Item(10,10).Value = 0.1133949940112261
$ws.Cells.Item(10,11).Value = 1.0
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.002632
$ws.Cells.Item(10,14).Value = 0.007896
$ws.Cells.Item(10,15).Value = 0.000006230194780080539
$ws.Cells.Item(10,16).Value = 0.000006230194780080539
$ws.Cells.Item(10,17).Value = 0.02557824888266667
$ws.Cells.Item(10,18).Value = 0.230204239944
$ws.Cells.Item(10,19).Value = 0.0000007064728997760049
$ws.Cells.Item(10,20).Value = 0.0000007064728997760049

# Row 11
$ws.Cells.Item(11,1).Value = "Neutro"
$ws.Cells.Item(11,2).Value = "Ccl12"
$ws.Cells.Item(11,3).Value = "Ccr1"
$ws.Cells.Item(11,4).Value = "M1"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 9.718179666666666
$ws.Cells.Item(11,8).Value = 29.154539
$ws.Cells.Item(11,9).Value = 0.1133949940112261
$ws.Cells.Item(11,10).Value = 0.1133949940112261
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 58.61543533333333
$ws.Cells.Item(11,14).Value = 175.846306
$ws.Cells.Item(11,15).Value = 0.1387483203821739
$ws.Cells.Item(11,16).Value = 0.1387483203821739
$ws.Cells.Item(11,17).Value = 569.6353318092148
$ws.Cells.Item(11,18).Value = 5126.717986282934
$ws.Cells.Item(11,19).Value = 0.01573336495880429
$ws.Cells.Item(11,20).Value = 0.01573336495880429

# Row 12
$ws.Cells.Item(12,1).Value = "Neutro"
$ws.Cells.Item(12,2).Value = "Ccl12"
$ws.Cells.Item(12,3).Value = "Ccr1"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 9.718179666666666
$ws.Cells.Item(12,8).Value = 29.154539
$ws.Cells.Item(12,9).Value = 0.1133949940112261
$ws.Cells.Item(12,10).Value = 0.1133949940112261
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 28.58678433333334
$ws.Cells.Item(12,14).Value = 85.76035300000001
$ws.Cells.Item(12,15).Value = 0.0676676422997042
$ws.Cells.Item(12,16).Value = 0.0676676422997042
$ws.Cells.Item(12,17).Value = 277.8115062435853
$ws.Cells.Item(12,18).Value = 2500.303556192267
$ws.Cells.Item(12,19).Value = 0.007673171893328748
$ws.Cells.Item(12,20).Value = 0.00767317189332875

# Row 13
$ws.Cells.Item(13,1).Value = "Neutro"
$ws.Cells.Item(13,2).Value = "Ccl12"
$ws.Cells.Item(13,3).Value = "Ccr1"
$ws.Cells.Item(13,4).Value = "Neutro"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 9.718179666666666
$ws.Cells.Item(13,8).Value = 29.154539
$ws.Cells.Item(13,9).Value = 0.1133949940112261
$ws.Cells.Item(13,10).Value = 0.1133949940112261
$ws.Cells.Item(13,11).Value = 3.0
$ws.Cells.Item(13,12).Value = 1.0
$ws.Cells.Item(13,13).Value = 335.2538503333333
$ws.Cells.Item(13,14).Value = 1005.761551
$ws.Cells.Item(13,15).Value = 0.7935778071233418
$ws.Cells.Item(13,16).Value = 0.7935778071233417
$ws.Cells.Item(13,17).Value = 3258.057151481109
$ws.Cells.Item(13,18).Value = 29322.51436332999
$ws.Cells.Item(13,19).Value = 0.08998775068619329
$ws.Cells.Item(13,20).Value = 0.08998775068619329

